$wb = $excel.ActiveWorkbook

# Update the two date values on the VoucherHeader sheet (DRAW DATE / PAY DATE)
$ws1 = $wb.Worksheets.Item("VoucherHeader")
$ws1.Range("A2").Value = "`n10/12/2020`n"
$ws1.Range("B2").Value = "10/12/2020`n"

# Move the selection on VoucherHeader to B3 (no longer the active tab)
$ws1.Range("B3").Select() | Out-Null

# Make lineTotals the active/selected sheet (tab) and keep its existing selection D28
$ws3 = $wb.Worksheets.Item("lineTotals")
$ws3.Select() | Out-Null
$ws3.Range("D28").Select() | Out-Null
